# fix orientation issues and app friendliness
#
# The sheet held a single subject's row of ventilation-summary data.
# Refresh row 2 with the next subject (Xe-020, re-processed 2021-11-19),
# and re-apply the text format (with its, until-now-implicit, border)
# across the header row and the leading ID/date columns so the row
# orientation / formatting stays consistent for downstream tooling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: subject identity + dates -------------------------------------
$ws.Range("A2").Value = "Xe-020"
$ws.Range("B2").Value = "2021-08-27"
$ws.Range("C2").Value = "2021-11-19"

# --- Row 2: numeric measurements ------------------------------------------
$ws.Range("D2").Value = 33.779210498172212
$ws.Range("E2").Value = 5.1195827471433963
$ws.Range("F2").Value = 12.14356792368665
$ws.Range("G2").Value = 1.5921490580928523
$ws.Range("H2").Value = 17.263150670830044
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 76.104038705692616
$ws.Range("K2").Value = 13.100916171979549
$ws.Range("L2").Value = 76.104038705692616
$ws.Range("M2").Value = 13.317091582884396
$ws.Range("N2").Value = 20.38225302817143
$ws.Range("O2").Value = 31.441512541605189
$ws.Range("P2").Value = 23.611158768829565
$ws.Range("Q2").Value = 8.5028994955907073
$ws.Range("R2").Value = 2.7450845829187114
$ws.Range("S2").Value = 81.70401125484679
$ws.Range("T2").Value = 4.5568404076450602
$ws.Range("U2").Value = 3.9014514634732183
$ws.Range("V2").Value = 3.9700785780461856
$ws.Range("W2").Value = 3.5514531791510824
$ws.Range("X2").Value = 2.3161651168376629
$ws.Range("Y2").Value = 5.5142815926139646
$ws.Range("Z2").Value = 13.571119446047316
$ws.Range("AA2").Value = 45.059145989613384
$ws.Range("AB2").Value = 35.855452971725335

# --- Re-apply text format/border across the header row and the ID/date
#     columns of row 2, mirroring the new, dedicated style used there ------
$ws.Range("A1:AB1").NumberFormat = "@"
$ws.Range("A2:C2").NumberFormat = "@"
